$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.752.70"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.07"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.23"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.689"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.08"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.348"
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.01"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0738"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0969"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.163.28"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.82"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.715"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.90"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.888.82"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "34.853.45"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.19"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "248.09"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.73"
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  -2.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -4.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.32"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.38"
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.23"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.36"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.25"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0581"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.55"
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.15"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.834"
$ws.Range("E38").Value = "  -8.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.97"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.31"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.88"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0662"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0210"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.08"
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.290.11"
$ws.Range("E45").Value = "  -5.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.36"
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.11"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0760"
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.45"
$ws.Range("E51").Value = "  -1.56%  "
